# Add the new data row (row 3) to the "AI Generated" sheet, as produced
# by the latest Katalon AI generation pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = "25 f4 Address"
$ws.Cells.Item(3, 2).Value = " "
$ws.Cells.Item(3, 3).Value = "25 f4 first"
$ws.Cells.Item(3, 4).Value = "25 f4 last"
$ws.Cells.Item(3, 5).Value = " "
$ws.Cells.Item(3, 6).Value = "25 f4 city"

# "2544" must be stored as text (matches the sibling cell G2, which holds
# the text "2502"), so force a text format before assigning, then drop the
# now-unneeded formatting so the cell keeps the sheet's default style.
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "2544"
$ws.Cells.Item(3, 7).ClearFormats()

$ws.Cells.Item(3, 8).Value = "{{address}}"
